$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Remove the BAU CCS / 45Q tax credit subsidy amount (policy removal)
$ws.Range("B11").Value = 0

# D3 on BCS-BCS carried a redundant explicit "General" number format;
# normalize it back to the plain Normal style (matches the other cells
# in the row/column which have no explicit style override).
$wsBcs = $wb.Worksheets.Item("BCS-BCS")
$wsBcs.Range("D3").Style = "Normal"

# Reflect the cell selection left behind after the edit
$ws.Activate()
$ws.Range("B12").Select()
